$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.466.03'
$ws.Range('E2').Value = '  -1.39%  '
$ws.Range('D3').Value = '2.184.73'
$ws.Range('E3').Value = '  -2.17%  '
$ws.Range('E4').Value = '  -0.09%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '251.89'
$ws.Range('E5').Value = '  +2.45%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.613'
$ws.Range('E6').Value = '  -0.93%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '75.21'
$ws.Range('E7').Value = '  -0.54%  '
$ws.Range('E8').Value = '  -0.01%  '
$ws.Range('E9').Value = '  -5.70%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '40.29'
$ws.Range('E10').Value = '  -2.03%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0909'
$ws.Range('E11').Value = '  -2.32%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.102'
$ws.Range('E12').Value = '  -0.03%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '6.77'
$ws.Range('E13').Value = '  -2.94%  '
$ws.Range('D14').Value = '2.511.70'
$ws.Range('E14').Value = '  -2.23%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '14.15'
$ws.Range('E15').Value = '  -3.81%  '
$ws.Range('D16').Value = '2.186.46'
$ws.Range('E16').Value = '  -2.01%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.769'
$ws.Range('E17').Value = '  -5.23%  '
$ws.Range('D18').Value = '42.376.53'
$ws.Range('E18').Value = '  -1.33%  '
$ws.Range('E19').Value = '  -3.28%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '70.69'
$ws.Range('E20').Value = '  -0.65%  '
$ws.Range('E21').Value = '  -2.10%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '226.85'
$ws.Range('E22').Value = '  -0.90%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '9.42'
$ws.Range('E23').Value = '  -9.90%  '
$ws.Range('E24').Value = '  -4.60%  '
$ws.Range('E25').Value = '  +0.00%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '10.44'
$ws.Range('E26').Value = '  -4.88%  '
$ws.Range('E27').Value = '  +0.55%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.27'
$ws.Range('E28').Value = '  +8.06%  '
$ws.Range('B29').Value = 'PancakeSwap'
$ws.Range('C29').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.16'
$ws.Range('E29').Value = '  -4.08%  '
$ws.Range('B30').Value = 'InjectiveProtocol'
$ws.Range('C30').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '37.57'
$ws.Range('E30').Value = '  +1.10%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '172.45'
$ws.Range('E31').Value = '  -1.07%  '
$ws.Range('E32').Value = '  -1.78%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0822'
$ws.Range('E33').Value = '  +3.68%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.13'
$ws.Range('E34').Value = '  -4.60%  '
$ws.Range('E35').Value = '  -1.78%  '
$ws.Range('E36').Value = '  -4.03%  '
$ws.Range('E37').Value = '  -3.57%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0336'
$ws.Range('E38').Value = '  +0.64%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '11.93'
$ws.Range('E39').Value = '  -8.71%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.06'
$ws.Range('E40').Value = '  -3.50%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.59'
$ws.Range('E41').Value = '  +11.42%  '
$ws.Range('E42').Value = '  -7.84%  '
$ws.Range('B43').Value = 'MultiversX'
$ws.Range('C43').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '58.74'
$ws.Range('E43').Value = '  -2.30%  '
$ws.Range('B44').Value = 'Algorand'
$ws.Range('C44').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.193'
$ws.Range('E44').Value = '  -3.38%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '101.26'
$ws.Range('E45').Value = '  -3.90%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0970'
$ws.Range('E46').Value = '  -2.24%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.460'
$ws.Range('E47').Value = '  +3.95%  '
$ws.Range('E48').Value = '  -4.71%  '
$ws.Range('E49').Value = '  -1.59%  '
$ws.Range('E50').Value = '  -2.50%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.65'
$ws.Range('E51').Value = '  -0.98%  '
